$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row just above the current row 272 (old rows 272-333 shift
# down to 273-334), then populate the newly-inserted row with the new weekly
# observation.
$ws.Rows.Item(272).Insert()

$ws.Cells.Item(272, 1).Value = 6
$ws.Cells.Item(272, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(272, 3).Value = "Metropolitana"
$ws.Cells.Item(272, 4).Value = 44543
$ws.Cells.Item(272, 5).Value = 13
$ws.Cells.Item(272, 6).Value = 100112039
$ws.Cells.Item(272, 7).Value = "Ciboulette"
$ws.Cells.Item(272, 8).Value = "Sin especificar"
$ws.Cells.Item(272, 9).Value = "Primera"
$ws.Cells.Item(272, 10).Value = 770
$ws.Cells.Item(272, 11).Value = 800
$ws.Cells.Item(272, 12).Value = 900
$ws.Cells.Item(272, 13).Value = 843
$ws.Cells.Item(272, 14).Value = "`$/docena de atados"
$ws.Cells.Item(272, 15).Value = "Región Metropolitana"
$ws.Cells.Item(272, 16).Value = 281
$ws.Cells.Item(272, 17).Value = 3
$ws.Cells.Item(272, 18).Value = "Hortaliza"
